$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2511.5908
$ws.Range("I138").Value = 1406.4286
$ws.Range("J138").Value = 4445.625
$ws.Range("K138").Value = 4219.2858
$ws.Range("L138").Value = 13336.875
$ws.Range("M138").Value = 920.7142000000003
$ws.Range("N138").Value = -23616.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11122.8
$ws.Range("I61").Value = 15666.667
$ws.Range("J61").Value = 4307
$ws.Range("K61").Value = 15666.667
$ws.Range("L61").Value = 4307
$ws.Range("M61").Value = -15454.667
$ws.Range("N61").Value = -4731

$ws.Range("H74").Value = 1272.4166
$ws.Range("I74").Value = 813.86664
$ws.Range("J74").Value = 2036.6666
$ws.Range("K74").Value = 813.86664
$ws.Range("L74").Value = 2036.6666
$ws.Range("M74").Value = 60.13336000000004
$ws.Range("N74").Value = -3784.6666

$ws.Range("H77").Value = 1272.4166
$ws.Range("I77").Value = 813.86664
$ws.Range("J77").Value = 2036.6666
$ws.Range("K77").Value = 4069.3332
$ws.Range("L77").Value = 10183.333
$ws.Range("M77").Value = 298.6668
$ws.Range("N77").Value = -18919.333

$ws.Range("H88").Value = 6192.222
$ws.Range("J88").Value = 6192.222
$ws.Range("L88").Value = 6192.222
$ws.Range("N88").Value = -7004.222

$ws.Range("H91").Value = 6192.222
$ws.Range("J91").Value = 6192.222
$ws.Range("L91").Value = 6192.222
$ws.Range("N91").Value = -9000.222

$ws.Range("H122").Value = 36701.473
$ws.Range("I122").Value = 36701.473
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 110104.419
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -107654.419
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 2235.5217
$ws.Range("I132").Value = 1866.5428
$ws.Range("J132").Value = 3409.5454
$ws.Range("K132").Value = 5599.6284
$ws.Range("L132").Value = 10228.6362
$ws.Range("M132").Value = -3069.6284
$ws.Range("N132").Value = -15288.6362

$ws.Range("H136").Value = 11122.8
$ws.Range("I136").Value = 15666.667
$ws.Range("J136").Value = 4307
$ws.Range("K136").Value = 47000.001
$ws.Range("L136").Value = 12921
$ws.Range("M136").Value = -44450.001
$ws.Range("N136").Value = -18021

$ws.Range("H139").Value = 78500
$ws.Range("J139").Value = 78500
$ws.Range("L139").Value = 78500
$ws.Range("N139").Value = -88780

$ws.Range("H140").Value = 97129.91
$ws.Range("J140").Value = 97129.91
$ws.Range("L140").Value = 97129.91
$ws.Range("N140").Value = -107489.91

$ws.Range("H141").Value = 59928.57
$ws.Range("J141").Value = 62615.383
$ws.Range("L141").Value = 62615.383
$ws.Range("N141").Value = -72975.383

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 10023.647
$ws.Range("I134").Value = 954
$ws.Range("J134").Value = 39500
$ws.Range("K134").Value = 2862
$ws.Range("L134").Value = 118500
$ws.Range("M134").Value = -327
$ws.Range("N134").Value = -123570

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 34225
$ws.Range("J95").Value = 34225
$ws.Range("L95").Value = 34225
$ws.Range("N95").Value = -39717

$ws.Range("H132").Value = 2219.1177
$ws.Range("I132").Value = 1540.6522
$ws.Range("J132").Value = 3637.7273
$ws.Range("K132").Value = 4621.9566
$ws.Range("L132").Value = 10913.1819
$ws.Range("M132").Value = -2091.9566
$ws.Range("N132").Value = -15973.1819

$ws.Range("H134").Value = 4181
$ws.Range("I134").Value = 3908
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 11724
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -9189
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6757602.5
$ws.Range("J131").Value = 8197577.5
$ws.Range("L131").Value = 24592732.5
$ws.Range("N131").Value = -24602812.5

$ws.Range("H141").Value = 2682
$ws.Range("I141").Value = 2503.6924
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 7511.0772
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -2331.0772
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 28000
$ws.Range("J64").Value = 28000
$ws.Range("L64").Value = 28000
$ws.Range("N64").Value = -28496

$ws.Range("H67").Value = 28000
$ws.Range("J67").Value = 28000
$ws.Range("L67").Value = 28000
$ws.Range("N67").Value = -29716

$ws.Range("H122").Value = 2402.9375
$ws.Range("I122").Value = 2265.1538
$ws.Range("K122").Value = 6795.4614
$ws.Range("M122").Value = -4345.4614

$ws.Range("H132").Value = 2206.513
$ws.Range("I132").Value = 1559.32
$ws.Range("J132").Value = 3362.2144
$ws.Range("K132").Value = 4677.96
$ws.Range("L132").Value = 10086.6432
$ws.Range("M132").Value = -2147.96
$ws.Range("N132").Value = -15146.6432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 45513.715
$ws.Range("J47").Value = 45513.715
$ws.Range("L47").Value = 45513.715
$ws.Range("N47").Value = -46493.715

$ws.Range("H52").Value = 45513.715
$ws.Range("J52").Value = 45513.715
$ws.Range("L52").Value = 45513.715
$ws.Range("N52").Value = -45979.715

$ws.Range("H132").Value = 4982.2354
$ws.Range("I132").Value = 4942.1333
$ws.Range("J132").Value = 5283
$ws.Range("K132").Value = 14826.3999
$ws.Range("L132").Value = 15849
$ws.Range("M132").Value = -12296.3999
$ws.Range("N132").Value = -20909

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 5083.3335
$ws.Range("I33").Value = 3750
$ws.Range("J33").Value = 5750
$ws.Range("K33").Value = 3750
$ws.Range("L33").Value = 5750
$ws.Range("M33").Value = -3500
$ws.Range("N33").Value = -6250

$ws.Range("H36").Value = 5083.3335
$ws.Range("I36").Value = 3750
$ws.Range("J36").Value = 5750
$ws.Range("K36").Value = 3750
$ws.Range("L36").Value = 5750
$ws.Range("M36").Value = -3500
$ws.Range("N36").Value = -6250

$ws.Range("H74").Value = 13885
$ws.Range("J74").Value = 13885
$ws.Range("L74").Value = 13885
$ws.Range("N74").Value = -15757

$ws.Range("H77").Value = 13885
$ws.Range("J77").Value = 13885
$ws.Range("L77").Value = 41655
$ws.Range("N77").Value = -51015

$ws.Range("H132").Value = 50002076
$ws.Range("I132").Value = 68183650
$ws.Range("J132").Value = 2751.625
$ws.Range("K132").Value = 204550950
$ws.Range("L132").Value = 8254.875
$ws.Range("M132").Value = -204548420
$ws.Range("N132").Value = -13314.875
